# Natmi following Dr Hou advice
# Rebuild the LR-pair table for Tgfa-Erbb2 (ECs/FAPs/sCs clusters), adding
# the sCs sending-cluster rows and refreshing the previously-computed
# specificity/expression statistics for the existing ECs rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A..T
$data = @(
    @(2, "ECs", "Tgfa", "Erbb2", "ECs", 3, 1, 1.570801666666667, 4.712405, 0.9976494164929963, 0.9976494164929962, 3, 1, 3.709791333333333, 11.129374, 0.4283284425582907, 0.4283284425582907, 5.827346409385557, 52.44611768447, 0.4273216207856326, 0.4273216207856326),
    @(3, "ECs", "Tgfa", "Erbb2", "FAPs", 3, 1, 1.570801666666667, 4.712405, 0.9976494164929963, 0.9976494164929962, 3, 1, 3.599001333333334, 10.797004, 0.4155367505499981, 0.4155367505499982, 5.653317292735557, 50.87985563462001, 0.4145599967176014, 0.4145599967176014),
    @(4, "ECs", "Tgfa", "Erbb2", "sCs", 3, 1, 1.570801666666667, 4.712405, 0.9976494164929963, 0.9976494164929962, 3, 1, 1.352297666666667, 4.056893, 0.1561348068917112, 0.1561348068917112, 2.124191428629445, 19.117722857665, 0.1557677989897623, 0.1557677989897623),
    @(5, "sCs", "Tgfa", "Erbb2", "ECs", 1, 0.3333333333333333, 0.003701, 0.011103, 0.002350583507003693, 0.002350583507003693, 3, 1, 3.709791333333333, 11.129374, 0.4283284425582907, 0.4283284425582907, 0.01372993772466667, 0.123569439522, 0.001006821772658097, 0.001006821772658097),
    @(6, "sCs", "Tgfa", "Erbb2", "FAPs", 1, 0.3333333333333333, 0.003701, 0.011103, 0.002350583507003693, 0.002350583507003693, 3, 1, 3.599001333333334, 10.797004, 0.4155367505499981, 0.4155367505499982, 0.01331990393466667, 0.119879135412, 0.0009767538323967333, 0.0009767538323967333),
    @(7, "sCs", "Tgfa", "Erbb2", "sCs", 1, 0.3333333333333333, 0.003701, 0.011103, 0.002350583507003693, 0.002350583507003693, 3, 1, 1.352297666666667, 4.056893, 0.1561348068917112, 0.1561348068917112, 0.005004853664333333, 0.04504368297899999, 0.0003670079019488629, 0.0003670079019488629)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value = $entry[$col]
    }
}
